$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing text storage so numeric-looking
# strings (e.g. "561.84") are not silently coerced into Number cells,
# matching the source data which stores every Price/Volume cell as text.
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$ws.Cells.Item(2, 4).Value = '59.172.08'
$ws.Cells.Item(2, 5).Value = '  +0.65%  '

$ws.Cells.Item(3, 4).Value = '2.977.06'
$ws.Cells.Item(3, 5).Value = '  -0.59%  '

$ws.Cells.Item(4, 5).Value = '  -0.14%  '

Set-TextValue $ws.Cells.Item(5, 4) '561.84'
$ws.Cells.Item(5, 5).Value = '  -0.03%  '

Set-TextValue $ws.Cells.Item(6, 4) '136.84'
$ws.Cells.Item(6, 5).Value = '  -0.27%  '

$ws.Cells.Item(7, 5).Value = '  -0.21%  '

Set-TextValue $ws.Cells.Item(8, 4) '0.516'
$ws.Cells.Item(8, 5).Value = '  -0.13%  '

$ws.Cells.Item(9, 4).Value = '2.971.46'
$ws.Cells.Item(9, 5).Value = '  -0.55%  '

$ws.Cells.Item(10, 5).Value = '  +0.70%  '

Set-TextValue $ws.Cells.Item(11, 4) '5.29'
$ws.Cells.Item(11, 5).Value = '  +8.77%  '

$ws.Cells.Item(12, 5).Value = '  -1.18%  '

$ws.Cells.Item(13, 5).Value = '  -0.51%  '

$ws.Cells.Item(14, 5).Value = '  -0.26%  '

$ws.Cells.Item(15, 5).Value = '  -0.39%  '

$ws.Cells.Item(16, 4).Value = '3.466.96'
$ws.Cells.Item(16, 5).Value = '  -0.74%  '

Set-TextValue $ws.Cells.Item(17, 4) '7.04'
$ws.Cells.Item(17, 5).Value = '  +0.63%  '

$ws.Cells.Item(18, 4).Value = '2.975.85'
$ws.Cells.Item(18, 5).Value = '  -0.81%  '

$ws.Cells.Item(19, 4).Value = '59.157.44'
$ws.Cells.Item(19, 5).Value = '  +0.44%  '

Set-TextValue $ws.Cells.Item(20, 4) '435.72'
$ws.Cells.Item(20, 5).Value = '  +2.15%  '

Set-TextValue $ws.Cells.Item(21, 4) '13.63'
$ws.Cells.Item(21, 5).Value = '  +0.13%  '

Set-TextValue $ws.Cells.Item(22, 4) '0.722'
$ws.Cells.Item(22, 5).Value = '  +1.03%  '

Set-TextValue $ws.Cells.Item(23, 4) '7.02'
$ws.Cells.Item(23, 5).Value = '  -1.81%  '

Set-TextValue $ws.Cells.Item(24, 4) '13.06'
$ws.Cells.Item(24, 5).Value = '  -3.00%  '

Set-TextValue $ws.Cells.Item(25, 4) '79.99'
$ws.Cells.Item(25, 5).Value = '  -0.38%  '

$ws.Cells.Item(26, 5).Value = '  +0.07%  '

$ws.Cells.Item(27, 5).Value = '  +5.37%  '

$ws.Cells.Item(28, 5).Value = '  -0.33%  '

$ws.Cells.Item(30, 5).Value = '  +0.91%  '

Set-TextValue $ws.Cells.Item(31, 4) '25.66'
$ws.Cells.Item(31, 5).Value = '  -0.32%  '

Set-TextValue $ws.Cells.Item(32, 4) '6.17'
$ws.Cells.Item(32, 5).Value = '  +1.19%  '

Set-TextValue $ws.Cells.Item(33, 4) '0.104'
$ws.Cells.Item(33, 5).Value = '  +5.73%  '

$ws.Cells.Item(34, 2).Value = 'PEPE'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(34, 4).Value = '0.0₃0778'
$ws.Cells.Item(34, 5).Value = '  +3.48%  '

$ws.Cells.Item(35, 2).Value = 'Mantle'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Cells.Item(35, 4) '0.994'
$ws.Cells.Item(35, 5).Value = '  +2.36%  '

$ws.Cells.Item(36, 2).Value = 'Filecoin'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Cells.Item(36, 4) '5.89'
$ws.Cells.Item(36, 5).Value = '  +2.13%  '

$ws.Cells.Item(37, 5).Value = '  -2.37%  '

Set-TextValue $ws.Cells.Item(38, 4) '48.53'
$ws.Cells.Item(38, 5).Value = '  -0.59%  '

Set-TextValue $ws.Cells.Item(39, 4) '8.70'
$ws.Cells.Item(39, 5).Value = '  -1.95%  '

Set-TextValue $ws.Cells.Item(40, 4) '2.79'
$ws.Cells.Item(40, 5).Value = '  +1.34%  '

Set-TextValue $ws.Cells.Item(41, 4) '395.10'
$ws.Cells.Item(41, 5).Value = '  +0.35%  '

$ws.Cells.Item(42, 5).Value = '  +0.11%  '

$ws.Cells.Item(43, 4).Value = '2.709.93'
$ws.Cells.Item(43, 5).Value = '  -0.74%  '

$ws.Cells.Item(44, 5).Value = '  -3.15%  '

$ws.Cells.Item(45, 5).Value = '  +1.15%  '

$ws.Cells.Item(46, 2).Value = 'Arweave'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue $ws.Cells.Item(46, 4) '35.32'
$ws.Cells.Item(46, 5).Value = '  +9.82%  '

$ws.Cells.Item(47, 2).Value = 'USDe'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Cells.Item(47, 4) '0.999'
$ws.Cells.Item(47, 5).Value = '  -0.02%  '

$ws.Cells.Item(48, 2).Value = 'Monero'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Cells.Item(48, 4) '122.05'
$ws.Cells.Item(48, 5).Value = '  -2.67%  '

$ws.Cells.Item(49, 5).Value = '  -0.25%  '

$ws.Cells.Item(50, 5).Value = '  -2.91%  '

Set-TextValue $ws.Cells.Item(51, 4) '23.09'
$ws.Cells.Item(51, 5).Value = '  -1.24%  '
